$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as Text so numeric-looking
# strings like "1.00" or "415.39" are not silently converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '66.624.71'
$ws.Range('D3').Value = '3.578.77'
$ws.Range('E3').Value = '  +3.09%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = '415.39'
$ws.Range('E5').Value = '  +0.27%  '
$ws.Range('D6').Value = '129.77'
$ws.Range('E6').Value = '  -0.30%  '
$ws.Range('D7').Value = '0.649'
$ws.Range('E7').Value = '  +3.61%  '
$ws.Range('D8').Value = '3.569.03'
$ws.Range('E8').Value = '  +2.93%  '
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('D10').Value = '0.776'
$ws.Range('E10').Value = '  +6.91%  '
$ws.Range('D11').Value = '0.175'
$ws.Range('E11').Value = '  +18.39%  '
$ws.Range('D12').Value = '0.0000332'
$ws.Range('E12').Value = '  +51.48%  '
$ws.Range('D13').Value = '42.36'
$ws.Range('E13').Value = '  -0.48%  '
$ws.Range('D14').Value = '9.89'
$ws.Range('E14').Value = '  +2.46%  '
$ws.Range('D15').Value = '4.127.44'
$ws.Range('E15').Value = '  +2.62%  '
$ws.Range('E16').Value = '  -0.36%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.646.35'
$ws.Range('E17').Value = '  +3.93%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').Value = '20.20'
$ws.Range('E18').Value = '  -1.51%  '
$ws.Range('D19').Value = '1.13'
$ws.Range('E19').Value = '  +4.24%  '
$ws.Range('D20').Value = '66.557.82'
$ws.Range('E20').Value = '  +6.15%  '
$ws.Range('D21').Value = '12.27'
$ws.Range('E21').Value = '  -3.47%  '
$ws.Range('D22').Value = '445.79'
$ws.Range('E22').Value = '  -4.34%  '
$ws.Range('D23').Value = '89.13'
$ws.Range('E23').Value = '  -2.13%  '
$ws.Range('D24').Value = '3.14'
$ws.Range('E24').Value = '  -3.51%  '
$ws.Range('D25').Value = '13.03'
$ws.Range('E25').Value = '  -1.82%  '
$ws.Range('D26').Value = '3.32'
$ws.Range('E26').Value = '  +0.73%  '
$ws.Range('D27').Value = '9.96'
$ws.Range('E27').Value = '  -5.33%  '
$ws.Range('D28').Value = '34.62'
$ws.Range('E28').Value = '  +3.66%  '
$ws.Range('D29').Value = '4.86'
$ws.Range('E29').Value = '  +1.25%  '
$ws.Range('E30').Value = '  +4.24%  '
$ws.Range('D31').Value = '12.34'
$ws.Range('E31').Value = '  +2.79%  '
$ws.Range('D32').Value = '0.117'
$ws.Range('E32').Value = '  +4.65%  '
$ws.Range('D33').Value = '7.32'
$ws.Range('E33').Value = '  -3.31%  '
$ws.Range('E34').Value = '  -4.44%  '
$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D35').Value = '39.72'
$ws.Range('E35').Value = '  -2.43%  '
$ws.Range('B36').Value = 'Dai'
$ws.Range('C36').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('D37').Value = '56.50'
$ws.Range('E37').Value = '  -3.37%  '
$ws.Range('D38').Value = '0.0490'
$ws.Range('E38').Value = '  +0.15%  '
$ws.Range('D39').Value = '0.0₃0729'
$ws.Range('E39').Value = '  +29.99%  '
$ws.Range('D40').Value = '0.147'
$ws.Range('E40').Value = '  +10.35%  '
$ws.Range('E41').Value = '  -0.33%  '
$ws.Range('D42').Value = '148.53'
$ws.Range('E42').Value = '  +1.46%  '
$ws.Range('D43').Value = '2.97'
$ws.Range('E43').Value = '  -3.77%  '
$ws.Range('D44').Value = '2.72'
$ws.Range('E44').Value = '  +0.83%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').Value = '4.31'
$ws.Range('E45').Value = '  -0.60%  '
$ws.Range('B46').Value = 'LidoDAOToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D46').Value = '3.23'
$ws.Range('E46').Value = '  -3.03%  '
$ws.Range('D47').Value = '0.309'
$ws.Range('E47').Value = '  -3.70%  '
$ws.Range('E48').Value = '  -5.25%  '
$ws.Range('D49').Value = '2.27'
$ws.Range('E49').Value = '  -4.97%  '
$ws.Range('D50').Value = '115.65'
$ws.Range('E50').Value = '  +5.54%  '
$ws.Range('D51').Value = '15.52'
$ws.Range('E51').Value = '  -4.87%  '

Write-Host "Applied cryptos update"
